# Auto-generated edit script: refresh the crypto price/volume snapshot
# (symbol-list update commit, Fri Feb 17 03:46:00 UTC 2023).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''309.16'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''-4.00%'
$ws.Range("E2").Style = "Normal"
# Row 3
$ws.Range("D3").Value = '''48.40'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''-6.35%'
$ws.Range("E3").Style = "Normal"
# Row 4
$ws.Range("D4").Value = '''5.165'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''-3.46%'
$ws.Range("E4").Style = "Normal"
# Row 5
$ws.Range("D5").Value = '''0.07762'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''-4.28%'
$ws.Range("E5").Style = "Normal"
# Row 6
$ws.Range("D6").Value = '''4.478'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''-2.00%'
$ws.Range("E6").Style = "Normal"
# Row 7
$ws.Range("D7").Value = '''1.315'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''18.55%'
$ws.Range("E7").Style = "Normal"
# Row 8
$ws.Range("D8").Value = '''1.555'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''-5.90%'
$ws.Range("E8").Style = "Normal"
# Row 9
$ws.Range("D9").Value = '''0.1228'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''-6.84%'
$ws.Range("E9").Style = "Normal"
# Row 10
$ws.Range("D10").Value = '''0.1947'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''-0.28%'
$ws.Range("E10").Style = "Normal"
# Row 11
$ws.Range("D11").Value = '''0.04683'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''2.95%'
$ws.Range("E11").Style = "Normal"
# Row 12
$ws.Range("D12").Value = '''0.09281'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''-3.39%'
$ws.Range("E12").Style = "Normal"
# Row 13
$ws.Range("D13").Value = '''0.1049'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''0.15%'
$ws.Range("E13").Style = "Normal"
# Row 14
$ws.Range("D14").Value = '''0.001266'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''-5.09%'
$ws.Range("E14").Style = "Normal"
# Row 15
$ws.Range("D15").Value = '''0.04162'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''-3.36%'
$ws.Range("E15").Style = "Normal"
# Row 16
$ws.Range("D16").Value = '''0.005858'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''0.97%'
$ws.Range("E16").Style = "Normal"
# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.329'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''-1.47%'
$ws.Range("E17").Style = "Normal"
# Row 18
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = '''2.276'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''-6.45%'
$ws.Range("E18").Style = "Normal"
# Row 19
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = '''0.3493'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''2.89%'
$ws.Range("E19").Style = "Normal"
# Row 20
$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D20").Value = '''8.278'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''1.03%'
$ws.Range("E20").Style = "Normal"
# Row 21
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D21").Value = '''0.1356'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''-2.99%'
$ws.Range("E21").Style = "Normal"
# Row 22
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D22").Value = '''0.3032'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''3.59%'
$ws.Range("E22").Style = "Normal"
# Row 23
$ws.Range("B23").Value = 'BitKan'
$ws.Range("C23").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D23").Value = '''0.001283'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''-1.69%'
$ws.Range("E23").Style = "Normal"
# Row 24
$ws.Range("B24").Value = 'HotbitToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D24").Value = '''0.004103'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''-3.79%'
$ws.Range("E24").Style = "Normal"
# Row 25
$ws.Range("D25").Value = '''0.0001349'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''0.13%'
$ws.Range("E25").Style = "Normal"
# Row 26
$ws.Range("E26").Value = '''-3.90%'
$ws.Range("E26").Style = "Normal"
# Row 38
$ws.Range("D38").Value = '''0.02573'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''-7.10%'
$ws.Range("E38").Style = "Normal"
# Row 39
$ws.Range("D39").Value = '''0.05876'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''6.10%'
$ws.Range("E39").Style = "Normal"
# Row 40
$ws.Range("D40").Value = '''0.01075'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''70.95%'
$ws.Range("E40").Style = "Normal"
# Row 41
$ws.Range("D41").Value = '''0.007915'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''1.84%'
$ws.Range("E41").Style = "Normal"
# Row 42
$ws.Range("E42").Value = '''-1.46%'
$ws.Range("E42").Style = "Normal"
# Row 43
$ws.Range("D43").Value = '''0.008425'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''9.70%'
$ws.Range("E43").Style = "Normal"
# Row 44
$ws.Range("D44").Value = '''0.007652'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''-13.18%'
$ws.Range("E44").Style = "Normal"
# Row 45
$ws.Range("D45").Value = '''0.3110'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''-11.84%'
$ws.Range("E45").Style = "Normal"
# Row 46
$ws.Range("D46").Value = '''0.00006887'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''1.00%'
$ws.Range("E46").Style = "Normal"
# Row 47
$ws.Range("E47").Value = '''0.11%'
$ws.Range("E47").Style = "Normal"
# Row 48
$ws.Range("D48").Value = '''0.05669'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''-6.15%'
$ws.Range("E48").Style = "Normal"
# Row 49
$ws.Range("E49").Value = '''0.29%'
$ws.Range("E49").Style = "Normal"
# Row 50
$ws.Range("D50").Value = '''0.00002099'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''0.11%'
$ws.Range("E50").Style = "Normal"
# Row 51
$ws.Range("D51").Value = '''0.0001999'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''0.11%'
$ws.Range("E51").Style = "Normal"
